$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Header text relabeling (sharedStrings content changes)
# -----------------------------------------------------------------
$ws.Range("A1").Value = "Roepnaam"
$ws.Range("B1").Value = "Tussenvoegsels"
$ws.Range("D1").Value = "E_Mail"
$ws.Range("E1").Value = "Persoon_ID"

# -----------------------------------------------------------------
# 2. Header row (row 1) formatting tweaks
#    - Wrap text on A1, B1, D1, E1 (not C1)
#    - Extend the light-green fill onto F1:F3
# -----------------------------------------------------------------
$ws.Range("A1").WrapText = $true
$ws.Range("B1").WrapText = $true
$ws.Range("D1").WrapText = $true
$ws.Range("E1").WrapText = $true

$ws.Range("F1:F3").Interior.Color = 14348002

# -----------------------------------------------------------------
# 3. Row 2 banding color swap: light-green -> white
# -----------------------------------------------------------------
$ws.Range("A2:E2").Interior.Color = 16777215

# -----------------------------------------------------------------
# 4. Row 4 picks up the explicit light-green banding fill
# -----------------------------------------------------------------
$ws.Range("A4:F4").Interior.Color = 14348002

# -----------------------------------------------------------------
# 5. Six new banded rows (5-10) appended below the existing table,
#    styled like the rest of the query-result block: light-green
#    fill plus a thin gray outline drawn around the whole new block.
# -----------------------------------------------------------------
$newRows = $ws.Range("A5:F10")
$newRows.Interior.Color = 14348002
$newRows.RowHeight = 13.55

foreach ($edge in 7,8,9,10) {
    $b = $newRows.Borders.Item($edge)
    $b.Color = 11184810
    $b.Weight = 2
    $b.LineStyle = 1
}

Write-Host "manegeplan export update applied"
